$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.033.31"
$ws.Range("E2").Value = "  +2.26%  "
$ws.Range("D3").Value = "2.589.47"
$ws.Range("E3").Value = "  +1.05%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "529.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.48"
$ws.Range("D6").Style = "Normal"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +1.07%  "
$ws.Range("D9").Value = "2.602.80"
$ws.Range("E9").Value = "  +1.11%  "
$ws.Range("E10").Value = "  +0.37%  "
$ws.Range("E11").Value = "  +2.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.333"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.94%  "
$ws.Range("E13").Value = "  +3.69%  "
$ws.Range("D14").Value = "3.049.17"
$ws.Range("E14").Value = "  +1.26%  "
$ws.Range("D15").Value = "58.989.84"
$ws.Range("E15").Value = "  +2.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.43"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.11%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000133"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.35%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "2.583.79"
$ws.Range("E18").Value = "  +0.70%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "346.57"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.24%  "
$ws.Range("E20").Value = "  +1.89%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.43"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.21%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "67.69"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.75%  "
$ws.Range("E25").Value = "  +0.99%  "
$ws.Range("E26").Value = "  +2.42%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.10"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.10%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("D30").Value = "0.0₃0722"
$ws.Range("E30").Value = "  +1.39%  "
$ws.Range("E31").Value = "  +3.96%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.88"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.57%  "
$ws.Range("E33").Value = "  +1.12%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "149.11"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.10%  "
$ws.Range("E35").Value = "  +1.37%  "
$ws.Range("E36").Value = "  +0.52%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "36.77"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.95%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.47"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.00%  "
$ws.Range("E39").Value = "  -0.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.816"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("E41").Value = "  +2.71%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.998"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.76"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.74%  "
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.595"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.71%  "
$ws.Range("B45").Value = "Bittensor"
$ws.Range("C45").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "268.58"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.36%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0957"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.80%  "
$ws.Range("E47").Value = "  +0.81%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.38"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.16%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.62"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.51%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "1.958.49"
$ws.Range("E50").Value = "  +0.64%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0221"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.39%  "
